# Automatische test-sync: 2025-06-24 19:42:50
$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 5 with the new e-mail log entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(5, 1).Value = "Herinnering betaling"
$logs.Cells.Item(5, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(5, 3).Value = "Ik zie dat ik nog een openstaande betaling heb. Kunt u dit bevestigen?"
$logs.Cells.Item(5, 4).Value = "Factuur / Administratie"
$logs.Cells.Item(5, 5).Value = "Beste klant,`nBedankt voor uw bericht. Om u beter van dienst te kunnen zijn, heb ik wat aanvullende informatie nodig. Kunt u mij uw klantnummer of factuurnummer doorgeven, zodat ik de openstaande betaling voor u kan controleren?`nAlvast bedankt voor uw medewerking.`nMet vriendelijke groet,`n[Naam] - E-mailassistent"
$logs.Cells.Item(5, 6).Value = "2025-06-24 19:42:28"
$logs.Cells.Item(5, 7).Value = "Ja"

# --- Sheet "Dashboard": re-sort category counts, Factuur / Administratie now has 2 ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(2, 1).Value = "Factuur / Administratie"
$dash.Cells.Item(2, 2).Value = 2

$dash.Cells.Item(3, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(3, 2).Value = 1

$dash.Cells.Item(4, 1).Value = "Sollicitatie / Vacature"
$dash.Cells.Item(4, 2).Value = 1

# --- Conditional formatting on "Logs" must now cover row 5 too ---
$catFcs = $logs.Range("D2:D4").FormatConditions
for ($i = 1; $i -le $catFcs.Count; $i++) {
    $catFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D5"))
}

$answeredFcs = $logs.Range("G2:G4").FormatConditions
for ($i = 1; $i -le $answeredFcs.Count; $i++) {
    $answeredFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G5"))
}
